# Update "想去人数" (want-to-go count) column F for a handful of events on
# the "展览" and "全部类型" sheets, reflecting freshly scraped counts.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (row -> new F value)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 244
$ws1.Range("F9").Value = 6778
$ws1.Range("F10").Value = 156
$ws1.Range("F16").Value = 16156
$ws1.Range("F20").Value = 178
$ws1.Range("F22").Value = 11330
$ws1.Range("F24").Value = 965
$ws1.Range("F25").Value = 4458
$ws1.Range("F26").Value = 309
$ws1.Range("F29").Value = 43

# Sheet "全部类型" (row -> new F value)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 244
$ws4.Range("F10").Value = 6778
$ws4.Range("F11").Value = 156
$ws4.Range("F18").Value = 16156
$ws4.Range("F22").Value = 178
$ws4.Range("F26").Value = 11330
$ws4.Range("F28").Value = 965
$ws4.Range("F29").Value = 4458
$ws4.Range("F30").Value = 309
$ws4.Range("F33").Value = 43
